# Edit script: apply revisions to Atsilo presentation
# 1) Add a new reviewer comment (Giulio) to slide 4 ("Tirocinanti" intro),
#    continuing the existing comment numbering (idx 11) right after idx 10.
# 2) On slide 15 ("Conclusioni"), merge three consecutive runs of the same
#    run-formatting into a single run (no visible text change).

$p = $ppt.ActivePresentation

# --- 1) New comment on slide 4 -------------------------------------------
$s4 = $p.Slides.Item(4)
# PowerPoint stores the comment marker position in points; use a value that
# resolves to the same near-origin marker ("10","10") used by the deck's
# other "general" (non-anchored) comments.
$commentLeft = 10 / 12700
$commentTop = 10 / 12700
$newComment = $s4.Comments.Add($commentLeft, $commentTop, "Giulio", "GF", "Ti consiglio di ispirarti alla slide di Marco, che ha messo proprio i requisiti funzionali presi dal RAD.`r`nQuesta parte qui ti conviene comunque metterla, perché spiega perché stai parlando di questo.")

# --- 2) Merge split runs on slide 15 --------------------------------------
$s15 = $p.Slides.Item(15)
$shape = $s15.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange
$mergedRange = $textRange.Characters(509, 372)
$mergedRange.Text = " traumatico, ma una volta che si ha preso confidenza con gli strumenti in nostro possesso poi è stato tutto più facile e rapido. Per le problematiche spiegate in precedenza e per la mancanza di tempo necessario, tutto il lavoro su l’individuazione degli attori, con annesso tutte le difficoltà sul capire bene i requisiti da adottare, poteva essere fatta in modo migliore."
